# SSO Login test data fixes
# - OPQA-5697 row (row 4): clientId in querystring updated, STORE value renamed
#   from "relayState" to "authenticationResponse.relayState" (now wrapped)
# - Row 10 / 11 (OPQA-5713 / OPQA-5714): the BODY JSON now references the new
#   "authenticationResponse.relayState" store key instead of "relayState"
# - Row 4 grows taller to fit the new wrapped STORE value

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 4 (OPQA-5697) ---------------------------------------------------
# QUERYSTRING: swap the clientId used to initiate login
$ws.Range("G4").Value = "?clientId=RAdneiT8SMGtuEnljBFWGA&app=dra&appurl=https://access.dev-stable.clarivate.com"

# STORE: rename the stored key, and make the cell wrap like its neighbours
$ws.Range("K4").Value = "authenticationResponse.relayState"
$ws.Range("K4").WrapText = $true

# The row now needs more vertical space to show the wrapped STORE text
$ws.Rows.Item(4).RowHeight = 60

# --- Rows 10 & 11 (OPQA-5713 / OPQA-5714) ---------------------------------
# BODY: the SAMLResponse payloads embed "(OPQA-5697_relayState)" as a
# placeholder that must track the STORE key renamed above.
$h10 = $ws.Range("H10").Text
$ws.Range("H10").Value = $h10.Replace("(OPQA-5697_relayState)", "(OPQA-5697_authenticationResponse.relayState)")

$h11 = $ws.Range("H11").Text
$ws.Range("H11").Value = $h11.Replace("(OPQA-5697_relayState)", "(OPQA-5697_authenticationResponse.relayState)")

# --- View state tidy-up ----------------------------------------------------
$ws.Activate()
$win = $ws.Application.ActiveWindow
$win.Split = $false
$win.ScrollColumn = 8
$win.ScrollRow = 1
$ws.Range("L2:L25").Select()
